$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 5-16 (players reshuffled), columns A=Name, B=Position, C=Team
$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
